$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5714677
$ws.Range("I19").Value = 15873376
$ws.Range("J19").Value = 408.75
$ws.Range("K19").Value = 15873376
$ws.Range("L19").Value = 408.75
$ws.Range("M19").Value = -15873201
$ws.Range("N19").Value = -758.75
$ws.Range("H33").Value = 505.5
$ws.Range("I33").Value = 329.26666
$ws.Range("K33").Value = 329.26666
$ws.Range("M33").Value = -100.26666
$ws.Range("H62").Value = 2300.889
$ws.Range("I62").Value = 1266.6666
$ws.Range("J62").Value = 2818
$ws.Range("K62").Value = 1266.6666
$ws.Range("L62").Value = 2818
$ws.Range("M62").Value = -642.6666
$ws.Range("N62").Value = -4066
$ws.Range("H65").Value = 2300.889
$ws.Range("I65").Value = 1266.6666
$ws.Range("J65").Value = 2818
$ws.Range("K65").Value = 6333.333000000001
$ws.Range("L65").Value = 14090
$ws.Range("M65").Value = -3213.333000000001
$ws.Range("N65").Value = -20330
$ws.Range("H132").Value = 6063420.5
$ws.Range("I132").Value = 6669280.5
$ws.Range("K132").Value = 20007841.5
$ws.Range("M132").Value = -20005311.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8334886.5
$ws.Range("I2").Value = 19232538
$ws.Range("J2").Value = 1388.1765
$ws.Range("K2").Value = 19232538
$ws.Range("L2").Value = 1388.1765
$ws.Range("M2").Value = -19232425
$ws.Range("N2").Value = -1614.1765
$ws.Range("H61").Value = 6466.4165
$ws.Range("I61").Value = 12403.8
$ws.Range("K61").Value = 12403.8
$ws.Range("M61").Value = -12191.8
$ws.Range("H88").Value = 2916.3845
$ws.Range("I88").Value = 1726.5
$ws.Range("J88").Value = 3445.2222
$ws.Range("K88").Value = 1726.5
$ws.Range("L88").Value = 3445.2222
$ws.Range("M88").Value = -1320.5
$ws.Range("N88").Value = -4257.2222
$ws.Range("H91").Value = 2916.3845
$ws.Range("I91").Value = 1726.5
$ws.Range("J91").Value = 3445.2222
$ws.Range("K91").Value = 1726.5
$ws.Range("L91").Value = 3445.2222
$ws.Range("M91").Value = -322.5
$ws.Range("N91").Value = -6253.2222
$ws.Range("H110").Value = 1849.95
$ws.Range("I110").Value = 764.64703
$ws.Range("K110").Value = 764.64703
$ws.Range("M110").Value = 1280.35297
$ws.Range("H116").Value = 8334886.5
$ws.Range("I116").Value = 19232538
$ws.Range("J116").Value = 1388.1765
$ws.Range("K116").Value = 19232538
$ws.Range("L116").Value = 1388.1765
$ws.Range("M116").Value = -19230244
$ws.Range("N116").Value = -5976.1765
$ws.Range("H122").Value = 1469.3864
$ws.Range("I122").Value = 1329.4073
$ws.Range("K122").Value = 3988.2219
$ws.Range("M122").Value = -1538.2219
$ws.Range("H132").Value = 2915.2744
$ws.Range("I132").Value = 2859.7896
$ws.Range("K132").Value = 8579.3688
$ws.Range("M132").Value = -6049.3688
$ws.Range("H136").Value = 6466.4165
$ws.Range("I136").Value = 12403.8
$ws.Range("K136").Value = 37211.39999999999
$ws.Range("M136").Value = -34661.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8334886.5
$ws.Range("I3").Value = 19232538
$ws.Range("J3").Value = 1388.1765
$ws.Range("K3").Value = 19232538
$ws.Range("L3").Value = 1388.1765
$ws.Range("M3").Value = -19232424
$ws.Range("N3").Value = -1616.1765

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1987.1482
$ws.Range("I31").Value = 1377
$ws.Range("K31").Value = 1377
$ws.Range("M31").Value = -1082
$ws.Range("H34").Value = 1987.1482
$ws.Range("I34").Value = 1377
$ws.Range("K34").Value = 1377
$ws.Range("M34").Value = -1175
$ws.Range("H58").Value = 13516902
$ws.Range("I58").Value = 2446.9524
$ws.Range("J58").Value = 31254624
$ws.Range("K58").Value = 2446.9524
$ws.Range("L58").Value = 31254624
$ws.Range("M58").Value = -2243.9524
$ws.Range("N58").Value = -31255030
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("H132").Value = 4136
$ws.Range("J132").Value = 2250
$ws.Range("L132").Value = 6750
$ws.Range("N132").Value = -11810
$ws.Range("H134").Value = 3990
$ws.Range("I134").Value = 4988.55
$ws.Range("K134").Value = 14965.65
$ws.Range("M134").Value = -12430.65
$ws.Range("H136").Value = 13516902
$ws.Range("I136").Value = 2446.9524
$ws.Range("J136").Value = 31254624
$ws.Range("K136").Value = 7340.8572
$ws.Range("L136").Value = 93763872
$ws.Range("M136").Value = -4790.8572
$ws.Range("N136").Value = -93768972

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 259.6842
$ws.Range("J98").Value = 429.6
$ws.Range("L98").Value = 1288.8
$ws.Range("N98").Value = -4284.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 38090.414
$ws.Range("I102").Value = 2891.75
$ws.Range("J102").Value = 51499.43
$ws.Range("K102").Value = 2891.75
$ws.Range("L102").Value = 51499.43
$ws.Range("M102").Value = -1269.75
$ws.Range("N102").Value = -54743.43
$ws.Range("H113").Value = 1722.84
$ws.Range("I113").Value = 1780.5
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 1780.5
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = 389.5
$ws.Range("N113").Value = -5640
$ws.Range("H126").Value = 297069.8
$ws.Range("I126").Value = 2050.3333
$ws.Range("K126").Value = 6150.999899999999
$ws.Range("M126").Value = -3680.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1968
$ws.Range("J61").Value = 4259.6
$ws.Range("L61").Value = 4259.6
$ws.Range("N61").Value = -4663.6
$ws.Range("H113").Value = 1968
$ws.Range("J113").Value = 4259.6
$ws.Range("L113").Value = 4259.6
$ws.Range("N113").Value = -8599.6
$ws.Range("H136").Value = 4206.1665
$ws.Range("I136").Value = 4610.091
$ws.Range("J136").Value = 3571.4285
$ws.Range("K136").Value = 13830.273
$ws.Range("L136").Value = 10714.2855
$ws.Range("M136").Value = -11280.273
$ws.Range("N136").Value = -15814.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3475.6875
$ws.Range("I136").Value = 3474
$ws.Range("K136").Value = 10422
$ws.Range("M136").Value = -7872
